$d = $word.ActiveDocument

# The title paragraph was accidentally typed as "Docker Composee" - the
# trailing "e" lives in its own run (a later edit appended it). Remove
# that stray trailing run/character so the title reads "Docker Compose"
# again, while leaving the original "Docker Compose" run untouched.
foreach ($para in $d.Paragraphs) {
    $text = $para.Range.Text
    if ($text -match "^Docker Composee\r?$") {
        $paraEnd = $para.Range.End
        $textEnd = $paraEnd - 1            # exclude the paragraph mark
        $lastChar = $d.Range($textEnd - 1, $textEnd)
        if ($lastChar.Text -eq "e") {
            $lastChar.Delete()
        }
        break
    }
}
